$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Price): force Text number format first so Excel does not
# auto-convert these dotted numeric-looking strings into real numbers
# (which would also silently strip meaningful trailing zeros).
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "26.163.02"
$ws.Range("D3").Value = "1.657.65"
$ws.Range("D5").Value = "215.42"
$ws.Range("D6").Value = "0.5251"
$ws.Range("D10").Value = "20.90"
$ws.Range("D11").Value = "0.07771"
$ws.Range("D12").Value = "1.662.23"
$ws.Range("D13").Value = "4.466"
$ws.Range("D14").Value = "1.882.43"
$ws.Range("D15").Value = "0.5530"
$ws.Range("D16").Value = "0.0₅8274"
$ws.Range("D17").Value = "65.21"
$ws.Range("D18").Value = "26.181.46"
$ws.Range("D20").Value = "4.758"
$ws.Range("D21").Value = "190.20"
$ws.Range("D22").Value = "10.30"
$ws.Range("D23").Value = "6.363"
$ws.Range("D24").Value = "1.002"
$ws.Range("D25").Value = "143.13"
$ws.Range("D26").Value = "0.1255"
$ws.Range("D27").Value = "7.420"
$ws.Range("D28").Value = "15.99"
$ws.Range("D29").Value = "1.429"
$ws.Range("D30").Value = "0.06147"
$ws.Range("D32").Value = "3.556"
$ws.Range("D35").Value = "0.9996"
$ws.Range("D36").Value = "2.400"
$ws.Range("D37").Value = "2.762"
$ws.Range("D38").Value = "0.5674"
$ws.Range("D40").Value = "5.918"
$ws.Range("D41").Value = "0.8539"
$ws.Range("D43").Value = "1.031.85"
$ws.Range("D44").Value = "99.66"
$ws.Range("D45").Value = "1.805.14"
$ws.Range("D47").Value = "56.13"
$ws.Range("D49").Value = "8.111"
$ws.Range("D51").Value = "5.982"

# Column E (Volume(1h)): percentage text, 2 leading/trailing spaces.
$ws.Range("E2").Value = "  +0.36%  "
$ws.Range("E3").Value = "  -0.11%  "
$ws.Range("E4").Value = "  -0.27%  "
$ws.Range("E5").Value = "  +3.95%  "
$ws.Range("E6").Value = "  +1.66%  "
$ws.Range("E7").Value = "  -0.25%  "
$ws.Range("E8").Value = "  +1.90%  "
$ws.Range("E9").Value = "  +1.57%  "
$ws.Range("E10").Value = "  +0.16%  "
$ws.Range("E11").Value = "  +3.01%  "
$ws.Range("E12").Value = "  +0.00%  "
$ws.Range("E13").Value = "  +1.71%  "
$ws.Range("E14").Value = "  -0.34%  "
$ws.Range("E15").Value = "  +2.92%  "
$ws.Range("E16").Value = "  +4.33%  "
$ws.Range("E17").Value = "  -1.36%  "
$ws.Range("E18").Value = "  +0.40%  "
$ws.Range("E19").Value = "  -0.22%  "
$ws.Range("E20").Value = "  +1.46%  "
$ws.Range("E21").Value = "  +1.66%  "
$ws.Range("E22").Value = "  +1.61%  "
$ws.Range("E23").Value = "  +2.92%  "
$ws.Range("E25").Value = "  -3.50%  "
$ws.Range("E26").Value = "  +3.79%  "
$ws.Range("E27").Value = "  +0.55%  "
$ws.Range("E28").Value = "  +2.36%  "
$ws.Range("E29").Value = "  +2.37%  "
$ws.Range("E30").Value = "  +2.83%  "
$ws.Range("E31").Value = "  +0.46%  "
$ws.Range("E32").Value = "  +2.53%  "
$ws.Range("E33").Value = "  +0.90%  "
$ws.Range("E34").Value = "  +1.90%  "
$ws.Range("E35").Value = "  +1.73%  "
$ws.Range("E36").Value = "  +0.42%  "
$ws.Range("E37").Value = "  +0.10%  "
$ws.Range("E38").Value = "  -3.36%  "
$ws.Range("E39").Value = "  +0.68%  "
$ws.Range("E40").Value = "  -0.61%  "
$ws.Range("E41").Value = "  +0.77%  "
$ws.Range("E43").Value = "  -6.36%  "
$ws.Range("E44").Value = "  -0.08%  "
$ws.Range("E45").Value = "  -0.64%  "
$ws.Range("E46").Value = "  -1.16%  "
$ws.Range("E47").Value = "  +2.03%  "
$ws.Range("E48").Value = "  +0.43%  "
$ws.Range("E49").Value = "  +1.32%  "
$ws.Range("E50").Value = "  -1.29%  "
$ws.Range("E51").Value = "  +2.28%  "

